$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor values updated
$ws.Range("B3").Value = 0.9809292892513893
$ws.Range("C3").Value = 0.9804085154786466
$ws.Range("D3").Value = 0.9810845970807667

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor, with new values
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9784095571646795
$ws.Range("C4").Value = 0.9778764950903678
$ws.Range("D4").Value = 0.9778514259950706

# Row 5: AdaBoostRegressor -> MLPRegressor, with new values
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.8349518389215028
$ws.Range("C5").Value = 0.846845853698444
$ws.Range("D5").Value = 0.853764028572104
